# Update odds values on Sheet1 for rows 4, 6, and 7 as per the FlashScore
# weekly update (2024-10-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4 ---
$ws.Range("I4").Value = 2.85
$ws.Range("K4").Value = 2.07
$ws.Range("L4").Value = 3.4
$ws.Range("M4").Value = 1.06
$ws.Range("N4").Value = 7.2
$ws.Range("P4").Value = 3.25
$ws.Range("Q4").Value = 1.9
$ws.Range("R4").Value = 1.85
$ws.Range("W4").Value = 8
$ws.Range("X4").Value = 11.75
$ws.Range("AB4").Value = 28
$ws.Range("AC4").Value = 7.2
$ws.Range("AD4").Value = 6.2
$ws.Range("AG4").Value = 9.75
$ws.Range("AH4").Value = 16
$ws.Range("AI4").Value = 10
$ws.Range("AJ4").Value = 37
$ws.Range("AK4").Value = 23
$ws.Range("AL4").Value = 29
$ws.Range("AP4").Value = 21
$ws.Range("AX4").Value = 15.5
$ws.Range("AY4").Value = 21
$ws.Range("AZ4").Value = 70
$ws.Range("BB4").Value = 250

# --- Row 6 ---
$ws.Range("G6").Value = 1.35
$ws.Range("H6").Value = 4.65
$ws.Range("I6").Value = 6.9
$ws.Range("J6").Value = 1.8
$ws.Range("K6").Value = 2.52
$ws.Range("L6").Value = 6.2
$ws.Range("O6").Value = 1.17
$ws.Range("P6").Value = 4.5
$ws.Range("Q6").Value = 1.52
$ws.Range("R6").Value = 2.37
$ws.Range("S6").Value = 1.28
$ws.Range("T6").Value = 3.35
$ws.Range("U6").Value = 1.78
$ws.Range("V6").Value = 1.93
$ws.Range("W6").Value = 8.5
$ws.Range("X6").Value = 7.4
$ws.Range("Z6").Value = 9.25
$ws.Range("AB6").Value = 22
$ws.Range("AD6").Value = 9.75
$ws.Range("AE6").Value = 18
$ws.Range("AF6").Value = 70
$ws.Range("AI6").Value = 22
$ws.Range("AK6").Value = 70
$ws.Range("AL6").Value = 55
$ws.Range("AM6").Value = 450
$ws.Range("AO6").Value = 6.1
$ws.Range("AQ6").Value = 16
$ws.Range("AT6").Value = 3.35
$ws.Range("AU6").Value = 7.9
$ws.Range("AV6").Value = 60
$ws.Range("AX6").Value = 37
$ws.Range("AY6").Value = 35
$ws.Range("BB6").Value = 400

# --- Row 7 ---
$ws.Range("G7").Value = 4.35
$ws.Range("H7").Value = 3.45
$ws.Range("I7").Value = 1.72
$ws.Range("J7").Value = 4.75
$ws.Range("K7").Value = 2.18
$ws.Range("L7").Value = 2.27
$ws.Range("M7").Value = 1.06
$ws.Range("N7").Value = 7.1
$ws.Range("O7").Value = 1.32
$ws.Range("P7").Value = 3.1
$ws.Range("Q7").Value = 1.95
$ws.Range("R7").Value = 1.78
$ws.Range("S7").Value = 1.39
$ws.Range("T7").Value = 2.77
$ws.Range("U7").Value = 1.88
$ws.Range("V7").Value = 1.83
$ws.Range("W7").Value = 11.5
$ws.Range("X7").Value = 24
$ws.Range("Y7").Value = 14.5
$ws.Range("Z7").Value = 75
$ws.Range("AA7").Value = 45
$ws.Range("AB7").Value = 50
$ws.Range("AC7").Value = 7.1
$ws.Range("AD7").Value = 6.9
$ws.Range("AE7").Value = 16.5
$ws.Range("AF7").Value = 80
$ws.Range("AG7").Value = 6.5
$ws.Range("AH7").Value = 7.8
$ws.Range("AI7").Value = 8.25
$ws.Range("AJ7").Value = 13.5
$ws.Range("AK7").Value = 14
$ws.Range("AL7").Value = 28
$ws.Range("AM7").Value = 700
$ws.Range("AN7").Value = 6.2
$ws.Range("AO7").Value = 25
$ws.Range("AP7").Value = 30
$ws.Range("AQ7").Value = 150
$ws.Range("AR7").Value = 175
$ws.Range("AS7").Value = 400
$ws.Range("AT7").Value = 2.77
$ws.Range("AU7").Value = 7.5
$ws.Range("AV7").Value = 70
$ws.Range("AW7").Value = 3.55
$ws.Range("AX7").Value = 8.5
$ws.Range("AY7").Value = 18
$ws.Range("AZ7").Value = 28
$ws.Range("BA7").Value = 60
